# Generate Report for Handoff
# Regenerates the localization-status report: the two previously "Ready for
# handoff" dependency (.png) rows collapse away and the run is re-keyed to
# a fresh pair of source .md files (new GUID-based names), with fresh
# handoff timestamps and xlf target names, on all three sheets.

$wb = $excel.ActiveWorkbook

$oldMd   = "70942ace-2c28-4d7e-9e24-d01c642723cb.md"
$newMd1  = "affb26ab-e351-47c8-b895-683175176dd9.md"
$newMd2  = "f8878dfd-f9fe-4b6f-8dde-91c59f7ccd04.md"

$zhHash1 = "affb26ab-e351-47c8-b895-683175176dd9.062b9bb8d6bf4a5ad978a9c6d66e826d675d57e4.zh-cn.xlf"
$zhHash2 = "f8878dfd-f9fe-4b6f-8dde-91c59f7ccd04.54240aecb9f3115c0a88325effd22a5f849288ff.zh-cn.xlf"
$deHash1 = "affb26ab-e351-47c8-b895-683175176dd9.062b9bb8d6bf4a5ad978a9c6d66e826d675d57e4.de-de.xlf"
$deHash2 = "f8878dfd-f9fe-4b6f-8dde-91c59f7ccd04.54240aecb9f3115c0a88325effd22a5f849288ff.de-de.xlf"

$zhTime = "2016-03-09 14:27:46"
$deTime = "2016-03-09 14:27:50"
$epoch  = "0001-01-01 00:00:00"

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/7dde6f719b8d3a4c838195f5623c75324ed8e6a6/e2e/"
$cfgUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/7dde6f719b8d3a4c838195f5623c75324ed8e6a6/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/934d353f300e98569d5b060819ff7c62ccce704e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5163cb42d5616d11b940579798b5d4a051156bd5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# Drop the old last row (.localization-config) since rows 3/4 (the two
# dependency pngs) collapse into a single replacement row, shrinking the
# sheet from 5 to 4 rows.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Rows("5:5").Delete()

$ov.Range("A2").Value = $newMd1
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

$ov.Range("A3").Value = $newMd2
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

$ov.Range("A4").Value = ".localization-config"
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), ($mdBase + $newMd1), "", "", $newMd1) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheets "zh-cn" / "de-de": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# Row 3 (the first dependency-png "IsDependency" row) is removed outright;
# the remaining rows shift up, and row 3 (old row 4) is converted from the
# dependency shape into a plain "Include" row for the second .md file,
# dropping its "Dependency From" (I) cell. Old row 5 (.localization-config)
# becomes row 4 unchanged in shape.
# ---------------------------------------------------------------------
function Update-LangSheet($sheetName, $xlfBase, $hash1, $hash2, $time) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Rows("3:3").Delete()

    $ws.Range("A2").Value = $newMd1
    $ws.Range("B2").Value = "Ready for handoff"
    $ws.Range("C2").Value = $hash1
    $ws.Range("D2").Value = $time
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = "Include"

    $ws.Range("A3").Value = $newMd2
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("C3").Value = $hash2
    $ws.Range("D3").Value = $time
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = "Include"
    $ws.Range("I3").ClearContents()

    # Row 4 (old row 5, the .localization-config row) already carries the
    # right text after the shift; nothing else to change there.

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), ($mdBase + $newMd1), "", "", $newMd1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), ($xlfBase + $hash1), "", "", $hash1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), ($mdBase + $newMd2), "", "", $newMd2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), ($xlfBase + $hash2), "", "", $hash2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null
}

Update-LangSheet "zh-cn" $zhXlfBase $zhHash1 $zhHash2 $zhTime
Update-LangSheet "de-de" $deXlfBase $deHash1 $deHash2 $deTime
